$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.153.90"
$ws.Range("E2").Value = "  -0.99%  "

# Row 3
$ws.Range("D3").Value = "2.915.93"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "371.66"
$ws.Range("E5").Value = "  +5.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.12"
$ws.Range("E6").Value = "  -2.45%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -2.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("E9").Value = "  -3.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.63"
$ws.Range("E10").Value = "  -2.60%  "

# Row 11
$ws.Range("E11").Value = "  +1.19%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0836"
$ws.Range("E12").Value = "  -1.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("E13").Value = "  -2.71%  "

# Row 14
$ws.Range("D14").Value = "3.370.53"
$ws.Range("E14").Value = "  -0.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.39"
$ws.Range("E15").Value = "  -1.94%  "

# Row 16
$ws.Range("D16").Value = "2.915.80"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.935"
$ws.Range("E17").Value = "  -3.35%  "

# Row 18
$ws.Range("D18").Value = "51.008.57"
$ws.Range("E18").Value = "  -1.26%  "

# Row 19
$ws.Range("E19").Value = "  -6.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  -1.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.89"
$ws.Range("E21").Value = "  -3.69%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0944"
$ws.Range("E22").Value = "  -1.79%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.26"
$ws.Range("E23").Value = "  -0.91%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.53"
$ws.Range("E24").Value = "  -0.72%  "

# Row 25
$ws.Range("E25").Value = "  -0.18%  "

# Row 26
$ws.Range("E26").Value = "  +3.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.172"
$ws.Range("E27").Value = "  -0.38%  "

# Row 28
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.73"
$ws.Range("E29").Value = "  -2.81%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("E30").Value = "  -9.22%  "

# Row 31
$ws.Range("E31").Value = "  -0.88%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.89"
$ws.Range("E33").Value = "  -3.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.13"
$ws.Range("E34").Value = "  -1.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "34.69"
$ws.Range("E35").Value = "  -2.31%  "

# Row 36
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0421"
$ws.Range("E38").Value = "  -0.78%  "

# Row 39
$ws.Range("E39").Value = "  -2.48%  "

# Row 40
$ws.Range("E40").Value = "  -0.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.08"
$ws.Range("E41").Value = "  -3.08%  "

# Row 42
$ws.Range("E42").Value = "  -5.48%  "

# Row 43
$ws.Range("E43").Value = "  -1.77%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.19"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.20"
$ws.Range("E45").Value = "  -0.17%  "

# Row 46
$ws.Range("E46").Value = "  -2.64%  "

# Row 47
$ws.Range("D47").Value = "2.019.37"
$ws.Range("E47").Value = "  -3.64%  "

# Row 48
$ws.Range("E48").Value = "  -0.63%  "

# Row 49
$ws.Range("E49").Value = "  -3.96%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.245"
$ws.Range("E50").Value = "  +3.14%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.201.98"
$ws.Range("E51").Value = "  -0.53%  "
